# Working on issues in TES.m
$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# Delete the now-unused sheets, leaving only "FF Code Base"
$wb.Worksheets.Item("TES Simulink Code").Delete() | Out-Null
$wb.Worksheets.Item("MATLAB App").Delete() | Out-Null

$ws = $wb.Worksheets.Item("FF Code Base")

# New row 4
$ws.Range("A4").Value = 44798
$ws.Range("B4").Value = "Outlet temperature increasing in time durring charging and holding modes"
$ws.Range("C4").Value = "g2 needs to be negative"
$ws.Range("D4").Value = 44798
$ws.Range("E4").Value = "TES.m, FF.m"

# New row 5
$ws.Range("A5").Value = 44798
$ws.Range("B5").Value = "Temp distribution goes unstable after a certain period of time"
$ws.Range("C5").Value = "Constant for BC1 contribution was being computed incorrectly in ""computeBCNow()"""
$ws.Range("D5").Value = 44798
$ws.Range("E5").Value = "TES.m "

# New row 6
$ws.Range("A6").Value = 44798
$ws.Range("B6").Value = "Droop characteristic not captured in discharge mode"
$ws.Range("E6").Value = "TES.m"

# Match the date-formatted style used by the existing rows (s="1", numFmtId 14)
$ws.Range("A2:A3").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null
$ws.Range("D2:D3").Copy() | Out-Null
$ws.Range("D4:D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("E9").Select() | Out-Null
